# Insert a new weekly record row at row 32 on the active sheet.
# This shifts the existing rows 32-47 down to 33-48 (preserving all of
# their data/formatting) and fills the newly-opened row 32 with the new
# "Inferno" / Ají record dated 2021-10-19 (serial 44488).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 32..47 down to 33..48, leaving row 32 empty (but formatted
# like the row above, matching Excel's native Insert behaviour).
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record.
$ws.Cells.Item(32, 1).Value = 11
$ws.Cells.Item(32, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(32, 3).Value = 'Bíobío'
$ws.Cells.Item(32, 4).Value = 44488
$ws.Cells.Item(32, 5).Value = 8
$ws.Cells.Item(32, 6).Value = 100112021
$ws.Cells.Item(32, 7).Value = 'Ají'
$ws.Cells.Item(32, 8).Value = 'Inferno'
$ws.Cells.Item(32, 9).Value = 'Primera'
$ws.Cells.Item(32, 10).Value = 50
$ws.Cells.Item(32, 11).Value = 27000
$ws.Cells.Item(32, 12).Value = 28000
$ws.Cells.Item(32, 13).Value = 27400
$ws.Cells.Item(32, 14).Value = '$/caja 12 kilos'
$ws.Cells.Item(32, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(32, 16).Value = 2283
$ws.Cells.Item(32, 17).Value = 12
$ws.Cells.Item(32, 18).Value = 'Hortaliza'
